# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rows 16-32 of Hoja1 hold one line per (worker, periodo mora). The
# "Valor Mora" (G) figures are refreshed for the 2306 period, and the
# rows for GISSETH DEL CARMEN DECHAMP MORALES / SANDY HELENA PUPO LEON /
# KEYRA LUZ NEGRETTE BAUTISTA are re-sorted so period 2306 for all three
# workers comes first (rows 24-26), followed by GISSETH's remaining
# periods 2307..2312 in ascending order (rows 27-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Salario Basico" (G) for the first block of workers (rows 16-23) ---
$ws.Range("G16").Value = 1000000
$ws.Range("G17").Value = 1000000
# G18 unchanged (1160000)
$ws.Range("G19").Value = 908526
$ws.Range("G20").Value = 1000000
$ws.Range("G21").Value = 908526
# G22 unchanged (1160000)
$ws.Range("G23").Value = 908526

# --- Re-sort rows 24-32: GISSETH / SANDY / KEYRA ---
# Row 24: GISSETH, periodo 2306 (was periodo 2312 here)
$ws.Range("C24").Value = "1102830646"
$ws.Range("D24").Value = "GISSETH DEL CARMEN DECHAMP MORALES"
$ws.Range("E24").Value = "2306"
$ws.Range("F24").Value = 46400
$ws.Range("G24").Value = 1000000

# Row 25: SANDY HELENA PUPO LEON, periodo 2306 (was GISSETH periodo 2311 here)
$ws.Range("C25").Value = "45564420"
$ws.Range("D25").Value = "SANDY HELENA PUPO LEON"
$ws.Range("E25").Value = "2306"
$ws.Range("F25").Value = 46400
$ws.Range("G25").Value = 908526

# Row 26: KEYRA LUZ NEGRETTE BAUTISTA, periodo 2306 (was GISSETH periodo 2310 here)
$ws.Range("C26").Value = "1047496345"
$ws.Range("D26").Value = "KEYRA LUZ NEGRETTE BAUTISTA"
$ws.Range("E26").Value = "2306"
$ws.Range("F26").Value = 46400
$ws.Range("G26").Value = 877803

# Row 27: GISSETH, periodo 2307 (F/G unchanged values, only E + G move)
$ws.Range("E27").Value = "2307"
$ws.Range("G27").Value = 1000000

# Row 28: GISSETH, periodo 2308
$ws.Range("E28").Value = "2308"
$ws.Range("G28").Value = 1000000

# Row 29: GISSETH, periodo 2309
$ws.Range("E29").Value = "2309"
$ws.Range("G29").Value = 1000000

# Row 30: GISSETH, periodo 2310
$ws.Range("E30").Value = "2310"
$ws.Range("G30").Value = 1000000

# Row 31: GISSETH, periodo 2311 (was SANDY periodo 2306 here)
$ws.Range("C31").Value = "1102830646"
$ws.Range("D31").Value = "GISSETH DEL CARMEN DECHAMP MORALES"
$ws.Range("E31").Value = "2311"
$ws.Range("F31").Value = 40000
$ws.Range("G31").Value = 1000000

# Row 32: GISSETH, periodo 2312 (was KEYRA periodo 2306 here)
$ws.Range("C32").Value = "1102830646"
$ws.Range("D32").Value = "GISSETH DEL CARMEN DECHAMP MORALES"
$ws.Range("E32").Value = "2312"
$ws.Range("F32").Value = 25333
$ws.Range("G32").Value = 1000000
